# Annual Family Gathering Registration - rework the header row.
# Old headers: A1 FIRST  NAME | B1 LAST NAME | C1 OTHER NAMES | D1 GENDER | E1 RESIDENCE | F1 CONTACT | G1 CHURCH
# New headers: A1 FULL NAME   | B1 CONTACT   | C1 DENOMINATION | D1 RESIDENCE   (columns E:G removed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three trailing columns (E:G) entirely - they are no longer used.
$ws.Range("E1:G1").EntireColumn.Delete()

# Overwrite the remaining four headers with the new wording / order.
$ws.Range("A1").Value = "FULL NAME"
$ws.Range("B1").Value = "CONTACT"
$ws.Range("C1").Value = "DENOMINATION"
$ws.Range("D1").Value = "RESIDENCE"

# Match the saved selection state from the authored workbook.
$ws.Range("D2").Select()
